$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(1.177551604473877,0,0.1111673006571401,0.1706298720381338,2.423515257680847,1.702612152364651,1.507559806245695,0,0.2655555384710722,0.5802664747389201,0.1919466336200912,0.2656339456856216,3.018603368224845),
  @(1.159121228193612,0,0.1112072717473254,0.1710572777252943,2.426490535055862,1.701411247377848,1.512107015662849,0,0.2660855716250907,0.5315427134649156,0.1858347664823725,0.2601866460158142,3.042393113258008),
  @(1.148327013197417,0,0.1112495207557416,0.1713375635083783,2.429343575623278,1.701507668663524,1.515470812905036,0,0.2664315478317461,0.5018435728428869,0.1821698538767507,0.2569619244733019,3.057787081099281),
  @(1.144060060659825,0,0.1112712077995113,0.1714562826130657,2.430764462194048,1.701756639637111,1.516985499993368,0,0.2665777115678054,0.4897959637178815,0.1806985391216926,0.255678100110682,3.064258194315034),
  @(1.143359509896982,0,0.1112750794332769,0.1714762679393593,2.431016002912166,1.701810645707766,1.517245708451767,0,0.2666022949005207,0.4877988016492623,0.1804555698891619,0.2554667540417519,3.065344674921707),
  @(1.148268933371526,0,0.1112497951099627,0.1713391463600056,2.429361692291273,1.701510177372441,1.515490657639035,0,0.2664334980745213,0.5016808712821614,0.182149921339871,0.2569444876300331,3.0578735512475),
  @(1.17108873119949,0,0.1111774182002918,0.170773541369212,2.42432825934722,1.702025034189845,1.509009084854981,0,0.2657340395046184,0.5634215868421393,0.189821090368639,0.2637308895118124,3.026642773138366),
  @(1.219963019112015,0,0.111175243739229,0.1698056273076491,2.422593579230266,1.709651867564332,1.500830794549344,0,0.2645247763319283,0.6862113036758899,0.2055579196290154,0.277986581061846,2.971640193286312),
  @(1.258366497067726,0,0.111257917202785,0.1691799832665495,2.426273286554817,1.719294743894181,1.497579831657163,0,0.2637345530402815,0.7774688693669987,0.2175401239452412,0.2890337909246838,2.935028618669186),
  @(1.276375422120395,0,0.1113136364834055,0.1689137932384488,2.429022093946372,1.724560075568206,1.496698687703827,0,0.2633962266279655,0.8192110829205319,0.2230819709683374,0.2941832136440894,2.91919672434102),
  @(1.283272022441878,0,0.1113373238533839,0.1688156326261605,2.430217432546513,1.726680338059424,1.496450876611789,0,0.263271140325358,0.8350504431195702,0.2251935534165881,0.296150900478267,2.913319882934147),
  @(1.281783298941463,0,0.1113321075289733,0.1688366560047456,2.429953129161547,1.726218079560397,1.496500430053175,0,0.263297945284366,0.8316377147974094,0.2247382094090682,0.2957263380273218,2.914580301691004),
  @(1.276941268990726,0,0.1113155334905898,0.1689056646508362,2.429117341080982,1.724731977483486,1.496676579931233,0,0.2633858750268327,0.8205135477365104,0.2232554320097506,0.294344742118092,2.918710860944302),
  @(1.273985400300461,0,0.1113057179180537,0.168948277956362,2.428625503212089,1.723838158525552,1.496795655408306,0,0.2634401288842074,0.8137038936120291,0.2223488789730936,0.293500777896277,2.921256361575502),
  @(1.257200374999087,0,0.1112546388254323,0.1691977492857235,2.426115262803606,1.71896832744882,1.49764943633032,0,0.2637570881358426,0.7747454855544333,0.2171797761233023,0.2886997490059926,2.93607983074908),
  @(1.247041001568959,0,0.1112279301594015,0.1693555030152827,2.424850514512798,1.716205933852351,1.498326229883702,0,0.2639569417076206,0.7509040555270587,0.2140319615601811,0.2857861467259895,2.945384326613762),
  @(1.241248360521439,0,0.1112142736877395,0.1694479729530751,2.424224246918513,1.7146997870653,1.498771770870974,0,0.2640738835852234,0.7372126445005165,0.2122300004650413,0.2841219974436129,2.950813479698297),
  @(1.239295805285593,0,0.1112099432963092,0.1694795797989845,2.424029585997133,1.714204035967327,1.498932290257898,0,0.2641138204493334,0.7325806810778488,0.2116213635953841,0.2835605543328654,2.95266499969491),
  @(1.248117234140949,0,0.11123059693535,0.1693385304463015,2.424974678188718,1.716491435053541,1.498248361292028,0,0.2639354609314353,0.7534397890849505,0.2143661646124428,0.2860950969187215,2.94438583087809),
  @(1.278361404436708,0,0.1113203315894289,0.1688853235655232,2.42935864238413,1.7251650510687,1.49662251100068,0,0.2633599657670143,0.82378010664155,0.2236906074502798,0.2947500705590045,2.917494403294143),
  @(1.298576369559981,0,0.111394053835312,0.16860450928659,2.433123857697197,1.731570507721784,1.496060332886429,0,0.2630015059445281,0.8699407626826599,0.2298604318263813,0.3005097820460065,2.900609167727687),
  @(1.287746373869822,0,0.1113533329132679,0.1687529805202923,2.431031984416052,1.728084371015427,1.496314621374268,0,0.2631912103857337,0.8452867713209855,0.2265605786285505,0.2974263129137071,2.909558004620671),
  @(1.247630519473972,0,0.1112293859940223,0.1693461982151847,2.424918229622889,1.716362104564496,1.498283389850798,0,0.2639451660248922,0.7522933361185551,0.2142150472363085,0.2859553865146367,2.94483700198662),
  @(1.206301617371025,0,0.1111609668253202,0.1700524177890474,2.422192448282516,1.706879634780989,1.502558597916121,0,0.2648346103718602,0.6528099306513013,0.2012266893481041,0.2740290205578262,2.985852331664503)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value2 = $row[$c]
    }
}

Write-Output "applied $($data.Count) rows"